$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title in A1
$ws.Range("A1").Value = "Iteration 2 Tomecards"

# Update Andrew Case's time spent value in B4
$ws.Range("B4").Value = "10h 30m"

# Set the active selection to B5, matching the saved view state
$ws.Range("B5").Select()
